# mensagens.xlsx: keep only the "Empresa"/"Mensagem" header plus the two
# GARCIA & MARQUES / GARCIA MADRUGA rows, refresh their send-date text from
# 06/02/25 -> 07/02/25, and widen column B to fit the long message text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original rows (besides header row 1):
#   2 ZANOTTOS BELEZA E TREINAMENTOS LTDA
#   3 PEDRADA FILMES LTDA
#   4 WR ASSESSORIA PECUARIA LTDA
#   5 PB SERVICOS DE FISIOTERAPIA LTDA
#   6 GARCIA & MARQUES M G LTDA      <- keep (becomes row 2)
#   7 GARCIA MADRUGA LTDA            <- keep (becomes row 3)
#   8 W M ALONSO LTDA
#   9 YOUNG TECHNOLOGY LTDA
# Delete the trailing rows first so earlier row numbers stay valid, then
# remove the leading block of rows we don't want to keep.
$ws.Rows("8:9").Delete()
$ws.Rows("2:5").Delete()

# Refresh the "valid as of" date mentioned inside each remaining message.
$b2 = $ws.Cells.Item(2, 2).Value()
$ws.Cells.Item(2, 2).Value = $b2.Replace("(06/02/25)", "(07/02/25)")

$b3 = $ws.Cells.Item(3, 2).Value()
$ws.Cells.Item(3, 2).Value = $b3.Replace("(06/02/25)", "(07/02/25)")

# Re-fit row heights (writing the long text can trigger an auto row-height
# bump) and widen column B so the long message text fits on one line.
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Range("B1:B3").EntireColumn.AutoFit()

# Match the saved selection from the source workbook.
$ws.Range("B8").Select()
